$wb = $excel.ActiveWorkbook

# Rename the first sheet from 2025_06_09 to 2025_06_10
$ws = $wb.Worksheets.Item(1)
$ws.Name = "2025_06_10"

# Select cell A2 on that sheet so the saved view records the selection
$ws.Activate()
$ws.Range("A2").Select()
